$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Ligand/Receptor/Edge expression metrics with recomputed TPM-based values
$ws.Range("I2").Value = 0.7783932765807232
$ws.Range("J2").Value = 0.7783932765807231
$ws.Range("M2").Value = 0.09551033333333332
$ws.Range("N2").Value = 0.286531
$ws.Range("O2").Value = 0.0198020999427218
$ws.Range("P2").Value = 0.0198020999427218
$ws.Range("Q2").Value = 0.02190281168133333
$ws.Range("R2").Value = 0.197125305132
$ws.Range("S2").Value = 0.01541382145759417
$ws.Range("T2").Value = 0.01541382145759417
$ws.Range("I3").Value = 0.7783932765807232
$ws.Range("J3").Value = 0.7783932765807231
$ws.Range("O3").Value = 0.07175622098770619
$ws.Range("P3").Value = 0.07175622098770619
$ws.Range("S3").Value = 0.05585455996967108
$ws.Range("T3").Value = 0.05585455996967107
$ws.Range("I4").Value = 0.7783932765807232
$ws.Range("J4").Value = 0.7783932765807231
$ws.Range("M4").Value = 4.381634666666667
$ws.Range("N4").Value = 13.144904
$ws.Range("O4").Value = 0.9084416790695721
$ws.Range("P4").Value = 0.9084416790695721
$ws.Range("Q4").Value = 1.004813988298667
$ws.Range("R4").Value = 9.043325894688001
$ws.Range("S4").Value = 0.707124895153458
$ws.Range("T4").Value = 0.7071248951534579
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.065288
$ws.Range("H5").Value = 0.195864
$ws.Range("I5").Value = 0.2216067234192769
$ws.Range("J5").Value = 0.2216067234192769
$ws.Range("M5").Value = 0.09551033333333332
$ws.Range("N5").Value = 0.286531
$ws.Range("O5").Value = 0.0198020999427218
$ws.Range("P5").Value = 0.0198020999427218
$ws.Range("Q5").Value = 0.006235678642666666
$ws.Range("R5").Value = 0.056121107784
$ws.Range("S5").Value = 0.004388278485127628
$ws.Range("T5").Value = 0.004388278485127627
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.065288
$ws.Range("H6").Value = 0.195864
$ws.Range("I6").Value = 0.2216067234192769
$ws.Range("J6").Value = 0.2216067234192769
$ws.Range("O6").Value = 0.07175622098770619
$ws.Range("P6").Value = 0.07175622098770619
$ws.Range("Q6").Value = 0.02259602446133333
$ws.Range("R6").Value = 0.203364220152
$ws.Range("S6").Value = 0.01590166101803512
$ws.Range("T6").Value = 0.01590166101803512
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.065288
$ws.Range("H7").Value = 0.195864
$ws.Range("I7").Value = 0.2216067234192769
$ws.Range("J7").Value = 0.2216067234192769
$ws.Range("M7").Value = 4.381634666666667
$ws.Range("N7").Value = 13.144904
$ws.Range("O7").Value = 0.9084416790695721
$ws.Range("P7").Value = 0.9084416790695721
$ws.Range("Q7").Value = 0.2860681641173333
$ws.Range("R7").Value = 2.574613477056
$ws.Range("S7").Value = 0.2013167839161142
$ws.Range("T7").Value = 0.2013167839161141